# Bug fix for one trial with wrong date
# Updates institution stats for trials affected by the corrected trial date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Aalborg University Hospital
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = 54.3
$ws.Range("E3").Value = 38.2
$ws.Range("F3").Value = 69.5

# Row 4: Aarhus University
$ws.Range("C4").Value = 115
$ws.Range("D4").Value = 60.2
$ws.Range("E4").Value = 53.1
$ws.Range("F4").Value = 66.90000000000001

# Row 36: Tampere University Hospital
$ws.Range("C36").Value = 17
$ws.Range("D36").Value = 70.8
$ws.Range("E36").Value = 50.8
$ws.Range("F36").Value = 85.09999999999999

# Row 44: University of Copenhagen
$ws.Range("C44").Value = 44
$ws.Range("D44").Value = 44.4
$ws.Range("E44").Value = 35
$ws.Range("F44").Value = 54.3

# Row 51: University of Tampere
$ws.Range("C51").Value = 6
$ws.Range("D51").Value = 66.7
$ws.Range("E51").Value = 35.4
$ws.Range("F51").Value = 87.90000000000001

# Row 56: Total
$ws.Range("C56").Value = 1097
$ws.Range("D56").Value = 51.9
$ws.Range("E56").Value = 49.8
$ws.Range("F56").Value = 54.1
